# Auto-generated edit script: updates cryptos list (prices & volumes),
# and swaps a few rows whose ranking order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.330.38'
$ws.Range("E2").Value = '  -0.31%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.841.38'
$ws.Range("E3").Value = '  -0.78%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.95'
$ws.Range("E5").Value = '  -0.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6277'
$ws.Range("E6").Value = '  -1.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07414'
$ws.Range("E8").Value = '  -2.31%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2897'
$ws.Range("E9").Value = '  -1.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.82'
$ws.Range("E10").Value = '  +0.80%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.837.95'
$ws.Range("E12").Value = '  -0.93%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.980'
$ws.Range("E13").Value = '  -1.02%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6773'
$ws.Range("E14").Value = '  -1.11%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001018'
$ws.Range("E15").Value = '  -2.72%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.99'
$ws.Range("E16").Value = '  -1.67%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.242'
$ws.Range("E17").Value = '  +1.27%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.355.85'
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.68'
$ws.Range("E19").Value = '  -0.88%  '

# Row 20
$ws.Range("E20").Value = '  -0.79%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9994'
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.412'
$ws.Range("E22").Value = '  -1.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.23'
$ws.Range("E24").Value = '  +0.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.469'
$ws.Range("E25").Value = '  -0.10%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1352'
$ws.Range("E26").Value = '  -3.56%  '

# Row 27
$ws.Range("E27").Value = '  -1.70%  '

# Row 28
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06491'
$ws.Range("E28").Value = '  +14.00%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.451'
$ws.Range("E29").Value = '  +2.00%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.485'
$ws.Range("E30").Value = '  +0.40%  '

# Row 31
$ws.Range("E31").Value = '  -2.17%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.069'
$ws.Range("E32").Value = '  -0.02%  '

# Row 33
$ws.Range("E33").Value = '  +0.34%  '

# Row 34
$ws.Range("E34").Value = '  -1.78%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6947'
$ws.Range("E35").Value = '  -1.07%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.566'
$ws.Range("E36").Value = '  -0.69%  '

# Row 37
$ws.Range("E37").Value = '  +1.49%  '

# Row 38
$ws.Range("E38").Value = '  +1.27%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.241.27'
$ws.Range("E39").Value = '  -0.93%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.736'
$ws.Range("E40").Value = '  +2.74%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9325'
$ws.Range("E41").Value = '  +3.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9991'
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.015.87'
$ws.Range("E43").Value = '  +0.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.81'
$ws.Range("E44").Value = '  -0.68%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.58'
$ws.Range("E45").Value = '  -0.76%  '

# Row 46
$ws.Range("E46").Value = '  +3.03%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.718'
$ws.Range("E47").Value = '  +2.22%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.049'
$ws.Range("E48").Value = '  -1.48%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.012'
$ws.Range("E49").Value = '  -0.53%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1151'
$ws.Range("E50").Value = '  -1.75%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3891'
$ws.Range("E51").Value = '  -2.03%  '
